$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the announcement text in B2
$ws.Range("B2").Value = "We will be having song practice on 9/23/2025 from 6:30pm - 8:30pm. This practice is specifically for the 40th year anniversary. We would like everyone to make it to this practice. Peb yuav muaj kawm nkauj rau hnub 9/23/2025 thaum 6:30pm txog 8:30pm. Qhov no yog special rau 40 xyoo. Xav kom sawvdaws tuaj."

# Increase row height of row 2 to fit the longer text
$ws.Rows(2).RowHeight = 90

# Move the active selection to B3
$ws.Range("B3").Select()
